# Update percent actual complete
# Adds a new "PercentActualComplete" worksheet at the end of the workbook
# with a small task table computing % of actual task completion.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new worksheet as the LAST sheet in the workbook.
# ---------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($sheetCount))
$ws.Name = "PercentActualComplete"

# Reference cells elsewhere in the workbook that already carry the exact
# cell styles we want to reuse, so the engine's format-dedup maps our new
# cells onto the *same* style indices instead of minting fresh ones.
$dateFmtSource = $wb.Worksheets.Item("Loc Phan").Range("K3")          # numFmtId 14 (mm-dd-yy), default font
$plainTnrSource = $wb.Worksheets.Item("Dang Nguyen").Range("B30")     # Times New Roman 12, no border/fill

# ---------------------------------------------------------------------------
# 2. Title / formula-description row.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Percentage of actual task complete per estimated`n"
$ws.Range("A1").HorizontalAlignment = 1
$ws.Rows("1:1").RowHeight = 15.75

$ws.Range("D1").Value = "(Actual Finish " + [char]0x2013 + " Start)( Plan Finish " + [char]0x2013 + " Start) /100%"
$plainTnrSource.Copy()
$ws.Range("D1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Header row (bold).
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Task"
$ws.Range("B3").Value = "Start"
$ws.Range("C3").Value = "Finish"
$ws.Range("D3").Value = "Actual Start"
$ws.Range("E3").Value = "Actual Finish"
$ws.Range("F3").Value = "% work complete"
$ws.Range("A3:F3").Font.Bold = $true

# ---------------------------------------------------------------------------
# 4. Data rows: Task name, Start, Finish, Actual Start, Actual Finish.
#    Column F is the computed "% work complete" formula.
# ---------------------------------------------------------------------------
$taskNames = @{
    4  = "Concept Operation"
    5  = "SRS"
    6  = "Architect Driver"
    7  = "Architect Design"
    8  = "Detail Design"
    9  = "System test cases"
    10 = "Integration test cases"
    11 = "Acceptance test"
    12 = "System test"
    13 = "Conduct System test"
    14 = "Conduct Acceptance test"
    15 = "Conduct Integration test"
    16 = "Programming"
}

$starts   = @{4=40817; 5=40847; 6=40848; 7=40870; 8=40890; 9=40899; 10=40910; 11=40883; 12=40899; 13=40983; 14=40994; 15=40983; 16=40875}
$finishes = @{4=40889; 5=40905; 6=40889; 7=40906; 8=40918; 9=41273; 10=40918; 11=40891; 12=40907; 13=40984; 14=40996; 15=40984; 16=40982}
$actStart = @{4=40927; 5=40858; 6=40851; 7=40928; 8=40978; 9=40950;           11=40932; 12=40950; 13=41030; 14=41034; 15=40974; 16=40978}
$actFinish= @{4=40929; 5=40925; 6=40928; 7=40950; 8=41016; 9=40954;           11=40935; 12=40954; 13=41034;           15=41030; 16=41027}
# row 10 has no Actual Start / Actual Finish at all; row 14 has an Actual
# Finish cell that is present but left empty (still date-formatted).

for ($r = 4; $r -le 16; $r++) {
    $ws.Range("A$r").Value = $taskNames[$r]

    $ws.Range("B$r").Value = $starts[$r]
    $plainTnrSource2 = $dateFmtSource
    $dateFmtSource.Copy()
    $ws.Range("B$r").PasteSpecial(-4122)

    $ws.Range("C$r").Value = $finishes[$r]
    $dateFmtSource.Copy()
    $ws.Range("C$r").PasteSpecial(-4122)

    if ($actStart.ContainsKey($r)) {
        $ws.Range("D$r").Value = $actStart[$r]
        $dateFmtSource.Copy()
        $ws.Range("D$r").PasteSpecial(-4122)
    }

    if ($r -eq 14) {
        # Actual Finish cell exists (date-formatted) but has no value.
        $dateFmtSource.Copy()
        $ws.Range("E$r").PasteSpecial(-4122)
    } elseif ($actFinish.ContainsKey($r)) {
        $ws.Range("E$r").Value = $actFinish[$r]
        $dateFmtSource.Copy()
        $ws.Range("E$r").PasteSpecial(-4122)
    }

    $ws.Range("F$r").FormulaR1C1 = "=((RC[-1]-RC[-4])*(RC[-3]-RC[-4]))/100"
}

Write-Host "PercentActualComplete sheet created"
